# Auto update stock data
$wb = $excel.ActiveWorkbook

# --- Norsk Hydro: refresh latest-row ratios ---
$wsNHY = $wb.Worksheets.Item("Norsk Hydro")
$wsNHY.Range("C2").Value = "'0.34"
$wsNHY.Range("D2").Value = "'5.06"
$wsNHY.Range("E2").Value = "'1.91"

# --- Reliance Steel & Aluminum: refresh EBITDA + backfill Altman Z-Score ---
$wsRS = $wb.Worksheets.Item("Reliance Steel & Aluminum")
$wsRS.Range("B2").Value = "'12.13"
$wsRS.Range("G2").Value = 5.92
$wsRS.Range("G3").Value = 5.92
$wsRS.Range("G4").Value = 5.92
$wsRS.Range("G5").Value = 5.92
$wsRS.Range("G6").Value = 5.92
$wsRS.Range("G7").Value = 5.92
$wsRS.Range("G8").Value = 5.92

# --- Kaiser Aluminum: refresh EBITDA ---
$wsKALU = $wb.Worksheets.Item("Kaiser Aluminum")
$wsKALU.Range("B2").Value = "'9.84"

# --- Ryerson Holding: refresh EBITDA ---
$wsRYI = $wb.Worksheets.Item("Ryerson Holding")
$wsRYI.Range("B2").Value = "'20.66"
